$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save current (pre-edit) values of rows 3, 4, 5 before overwriting,
# since the edit cyclically rotates the data among these rows:
#   new Row3 = old Row5
#   new Row4 = old Row3
#   new Row5 = old Row4

$oldRow3 = @{
    D = $ws.Range("D3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    N = $ws.Range("N3").Value2
    O = $ws.Range("O3").Value2
    P = $ws.Range("P3").Value2
    S = $ws.Range("S3").Value2
}

$oldRow4 = @{
    D = $ws.Range("D4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    N = $ws.Range("N4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
    S = $ws.Range("S4").Value2
}

$oldRow5 = @{
    D = $ws.Range("D5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    N = $ws.Range("N5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
    S = $ws.Range("S5").Value2
}

# New Row 3 <- old Row 5
$ws.Range("D3").Value = $oldRow5.D
$ws.Range("L3").Value = $oldRow5.L
$ws.Range("M3").Value = $oldRow5.M
$ws.Range("N3").Value = $oldRow5.N
$ws.Range("O3").Value = $oldRow5.O
$ws.Range("P3").Value = $oldRow5.P
$ws.Range("S3").Value = $oldRow5.S

# New Row 4 <- old Row 3
$ws.Range("D4").Value = $oldRow3.D
$ws.Range("L4").Value = $oldRow3.L
$ws.Range("M4").Value = $oldRow3.M
$ws.Range("N4").Value = $oldRow3.N
$ws.Range("O4").Value = $oldRow3.O
$ws.Range("P4").Value = $oldRow3.P
$ws.Range("S4").Value = $oldRow3.S

# New Row 5 <- old Row 4
$ws.Range("D5").Value = $oldRow4.D
$ws.Range("L5").Value = $oldRow4.L
$ws.Range("M5").Value = $oldRow4.M
$ws.Range("N5").Value = $oldRow4.N
$ws.Range("O5").Value = $oldRow4.O
$ws.Range("P5").Value = $oldRow4.P
$ws.Range("S5").Value = $oldRow4.S
